$wb = $excel.ActiveWorkbook

# --- Sheet 1: "Weekly Timesheet" ---
$ws1 = $wb.Worksheets.Item("Weekly Timesheet")

# Fix employee last names (Client column, B)
$ws1.Range("B3").Value = "Evans"
$ws1.Range("B5").Value = "Hewett"
$ws1.Range("B6").Value = "Howard"

# Simulator full-month coverage: populate Rate (E) and Total (F) for each day row
$ws1.Range("E2:E6").Value = 140
$ws1.Range("F2:F6").Value = 1120

# Persist computed subtotal/grand-total figures
$ws1.Range("F8").Value = 5600
$ws1.Range("F12").Value = 5600
$ws1.Range("F13").Value = 5600

# --- Sheet 2: "Jason Schema" ---
$ws2 = $wb.Worksheets.Item("Jason Schema")

$ws2.Range("F2:F6").Value = 140
$ws2.Range("G2:G6").Value = 1120

# --- Fix employee id (appears once per day row on the Jason Schema sheet) ---
$ws2.Range("B2:B6").Value = "emp_jp4mlvog"
